# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-146) on the sheet: 2023-10-25 (serial 45224) -> 2023-11-03 (serial 45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 146; $r++) {
    $ws.Cells.Item($r, 3).Value = 45233
}
